# Actualización automática 2025-08-25 17:20:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D13").Value = 1831.68
$ws1.Range("M22").Value = 650.9299999999999
$ws1.Range("D27").Value = 9331.200000000001
$ws1.Range("D34").Value = "6 de 32"
$ws1.Range("M34").Value = "6 de 32"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 1831.68
$ws2.Range("F22").Value = 650.9299999999999
$ws2.Range("F34").Value = 30741.96

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 16348.61
$ws3.Range("E3").Value = -13228.4955
$ws3.Range("F3").Value = 5.239746810573779

$ws3.Range("D16").Value = 7933.63
$ws3.Range("E16").Value = 13939.47
$ws3.Range("F16").Value = 0.3627117326762096

$ws3.Range("D19").Value = 31260.17000000001
$ws3.Range("E19").Value = 849.1110755578718
$ws3.Range("F19").Value = 0.9735555874465148
